$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.99%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'-1.30%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.102"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.38%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07959"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'5.70%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.306"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.08%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.768"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.13%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-0.49%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9224"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.73%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1734"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.35%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07544"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.43%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09322"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'12.59%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03044"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.40%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1002"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.82%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'0.48%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005758"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.75%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.484"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.46%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'1.39%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.22%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1335"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.21%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'3.899"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-16.19%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1700"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'8.81%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04629"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.22%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001251"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.80%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004479"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.34%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-7.58%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003395"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'24.11%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01762"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-1.19%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04633"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.10%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006970"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-4.58%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1362"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.15%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002188"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.34%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-4.68%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006264"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.66%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.007973"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-19.28%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.156"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'40.94%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E50").Style = "Normal"
